$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Predicted Date of Failure column (G) stores its dates as plain text,
# matching the same cell style already used by the Reference Date column (E).
# Temporarily force text format so Excel doesn't auto-convert the string to
# a date serial, then clear the formatting again so no stray style gets
# attached to the cells (they keep no explicit style in the original file).
$ws.Range("G2:G5").NumberFormat = "@"

# Update Defect Size (mm), Number of Days before Failure,
# Predicted Date of Failure, and Comment for rows 2-5.
$ws.Range("D2").Value = 30.1317
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "2017-04-07"
$ws.Range("H2").Value = "The Flange Height has violated the wheel alarm settings thresholds"

$ws.Range("D3").Value = 30.064
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "2017-04-07"
$ws.Range("H3").Value = "The Flange Height has violated the wheel alarm settings thresholds"

$ws.Range("D4").Value = 30.2994
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "2017-04-07"
$ws.Range("H4").Value = "The Flange Height has violated the wheel alarm settings thresholds"

$ws.Range("D5").Value = 30.0098
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "2017-04-07"
$ws.Range("H5").Value = "The Flange Height has violated the wheel alarm settings thresholds"

$ws.Range("G2:G5").ClearFormats()
